$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D:E (rows 2-51) to Text format temporarily so that
# numeric-looking strings (e.g. "29.166.56", "0.9989", "1.850") are
# preserved exactly as text, matching the source data which stores
# these as inline/shared strings, not numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.166.56'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.840.42'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '241.25'
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").Value = '0.6870'
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.3019'
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = '0.07488'
$ws.Range("D10").Value = '23.22'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").Value = '0.07664'
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").Value = '1.832.89'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '5.069'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("D14").Value = '0.6840'
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '87.95'
$ws.Range("E15").Value = '  -5.89%  '
$ws.Range("D16").Value = '6.164'
$ws.Range("E16").Value = '  -7.14%  '
$ws.Range("D17").Value = '29.152.92'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '0.000008181'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").Value = '2.081.26'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = '228.26'
$ws.Range("E20").Value = '  -5.48%  '
$ws.Range("E21").Value = '  -1.78%  '
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '7.426'
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '0.1458'
$ws.Range("E25").Value = '  -4.01%  '
$ws.Range("D26").Value = '160.06'
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").Value = '8.774'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("D29").Value = '1.517'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").Value = '4.279'
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("D31").Value = '4.158'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").Value = '1.196'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '0.05188'
$ws.Range("E33").Value = '  +1.30%  '
$ws.Range("D34").Value = '0.7671'
$ws.Range("E34").Value = '  -3.19%  '
$ws.Range("D35").Value = '1.850'
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").Value = '1.136'
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("D37").Value = '2.673'
$ws.Range("D38").Value = '1.318.17'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").Value = '0.01838'
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("D40").Value = '2.722'
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").Value = '0.9347'
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '104.83'
$ws.Range("E42").Value = '  -2.69%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.778'
$ws.Range("E43").Value = '  -4.58%  '
$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '0.00000000124'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").Value = '65.34'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").Value = '1.982.65'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Value = '0.5196'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = '9.556'
$ws.Range("E49").Value = '  -2.07%  '
$ws.Range("D50").Value = '1.775'
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").Value = '0.05941'
$ws.Range("E51").Value = '  +1.00%  '

# Restore original (default/General) formatting so no stray cell styles
# are introduced compared to the source workbook.
$dataRange.ClearFormats()
